$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B35").Value = "Sie werden eine Liste mit Name sehen und sollen dort ankreuzen, welche davon **Komponisten/Komponistinnen** des 19. Jhdt. sind (Hauptwirkungzeit). Sie haben dazu {{time_out}} Sekunden Zeit."
$ws.Range("B36").Value = "Bitte wählen Sie alle **Komponisten/Komponistinnen** des 19. Jhdt. (Hauptwirkungzeit) aus der untenstehenden Liste aus.  Sie haben {{time_out}} Sekunden Zeit."

$ws.Rows.Item(35).RowHeight = 90

$ws.Range("B36").Select()
